# "Medications quick change per Brandie H., LOCUS WIP part 2"
#
# Adds four new UM description rows (S0280/S0315/S0316/S9470) to the bottom
# of the "desc" sheet, restores the hidden _FilterDatabase name that Excel
# writes out whenever AutoFilter has been used on the sheet, and leaves the
# selection where the author ended up after the edit (row ~96-105 area).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new rows (columns B and C only, column A left blank like
#     the other CPT-code rows further up the sheet) ---------------------
$ws.Range("B102").Value = "S0280"
$ws.Range("C102").Value = "Medical Home Prog"

$ws.Range("B103").Value = "S0315"
$ws.Range("C103").Value = "Disease Mgmt Prog"

$ws.Range("B104").Value = "S0316"
$ws.Range("C104").Value = "Disease Mgmt Prog"

$ws.Range("B105").Value = "S9470"
$ws.Range("C105").Value = "Health | Dietician Svcs"

# --- Recreate the hidden AutoFilter-backed defined name over the
#     original data extent (desc!$A$1:$C$101) ---------------------------
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=desc!`$A`$1:`$C`$101")
$filterName.Visible = $false

# --- Leave the active selection where the author ended up working ------
$ws.Range("C96").Select()

$wb.Save()
